$d = $word.ActiveDocument

# 1. Append a new sentence (as its own run) after
#    "...output mat files as necessary for your study"
$findRange = $d.Content
$findRange.Find.Execute("output mat files as necessary for your study", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$insertPointVpn = $d.Range($findRange.End, $findRange.End)
$insertPointVpn.InsertAfter(" (Ensure you are connected to VPN and have X Drive with Aqua landmarks for data analysis)")

# 2. Remove the "(To be incorporated) " prefix from the Plot Event Raster bullet
$d.Content.Find.Execute("(To be incorporated) Plot Event Raster", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Plot Event Raster", 2)

# 3. Split the single "Plot Event Raster: ..." run into two separate runs:
#    "Plot " and "Event Raster: ...". Clear the paragraph's text and
#    re-insert the two pieces individually so each becomes its own <w:r>.
$p = $d.Paragraphs.Last
$full = $p.Range
$wholeRange = $d.Range($full.Start, $full.End - 1)
$wholeRange.Text = ""

$part1 = "Plot "
$part2 = "Event Raster: plots multiple figures that analyzes event wise and stimulus wise data such as size, orientation, distributions and more "

$insertPoint1 = $d.Range($full.Start, $full.Start)
$insertPoint1.InsertAfter($part1)

$insertPoint2 = $d.Range($full.Start + $part1.Length, $full.Start + $part1.Length)
$insertPoint2.InsertAfter($part2)

# 4. Remove the stray _GoBack bookmark left over from the prior edit location
$d.Bookmarks("_GoBack").Delete()
